# Fix for ical file.
# Remove the "NO CLASS - Happy Thanksgiving" row from the Meetups sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetups")

# Row 15 currently holds: 44524 | (blank) | (blank) | NO CLASS - Happy Thanksgiving
# Delete the entire row, shifting the rows below it up.
$ws.Rows.Item(15).Delete()

# Update the selection to match the author's resulting cursor position.
$ws.Activate()
$ws.Range("B24").Select()
